$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 92
$ws_ALC.Range("H92").Value = 4902555.5
$ws_ALC.Range("I92").Value = 6944855
$ws_ALC.Range("K92").Value = 6944855
$ws_ALC.Range("M92").Value = -6943607

# ALC row 100
$ws_ALC.Range("H100").Value = 2946.923
$ws_ALC.Range("I100").Value = 2385
$ws_ALC.Range("J100").Value = 3428.5715
$ws_ALC.Range("K100").Value = 2385
$ws_ALC.Range("L100").Value = 3428.5715
$ws_ALC.Range("M100").Value = -1844
$ws_ALC.Range("N100").Value = -4510.5715

# ALC row 121
$ws_ALC.Range("H121").Value = 2133.5715
$ws_ALC.Range("J121").Value = 1640
$ws_ALC.Range("L121").Value = 4920
$ws_ALC.Range("N121").Value = -8414

# ALC row 132
$ws_ALC.Range("H132").Value = 5516.8613
$ws_ALC.Range("I132").Value = 5946.6787
$ws_ALC.Range("J132").Value = 4012.5
$ws_ALC.Range("K132").Value = 17840.0361
$ws_ALC.Range("L132").Value = 12037.5
$ws_ALC.Range("M132").Value = -15310.0361
$ws_ALC.Range("N132").Value = -17097.5

# ALC row 135
$ws_ALC.Range("H135").Value = 928.0732
$ws_ALC.Range("I135").Value = 652.4838999999999
$ws_ALC.Range("J135").Value = 1782.4
$ws_ALC.Range("K135").Value = 5872.3551
$ws_ALC.Range("L135").Value = 16041.6
$ws_ALC.Range("M135").Value = -3337.3551
$ws_ALC.Range("N135").Value = -21111.6

# ALC row 137
$ws_ALC.Range("H137").Value = 41668052
$ws_ALC.Range("I137").Value = 28572494
$ws_ALC.Range("J137").Value = 76925330
$ws_ALC.Range("K137").Value = 85717482
$ws_ALC.Range("L137").Value = 230775990
$ws_ALC.Range("M137").Value = -85714932
$ws_ALC.Range("N137").Value = -230781090

# ALC row 138
$ws_ALC.Range("H138").Value = 1960.3627
$ws_ALC.Range("I138").Value = 606.55
$ws_ALC.Range("J138").Value = 3022.1765
$ws_ALC.Range("K138").Value = 1819.65
$ws_ALC.Range("L138").Value = 9066.529500000001
$ws_ALC.Range("M138").Value = 3320.35
$ws_ALC.Range("N138").Value = -19346.5295

# ARM row 32
$ws_ARM.Range("H32").Value = 10697.506
$ws_ARM.Range("I32").Value = 3816.6826
$ws_ARM.Range("J32").Value = 26179.357
$ws_ARM.Range("K32").Value = 3816.6826
$ws_ARM.Range("L32").Value = 26179.357
$ws_ARM.Range("M32").Value = -3529.6826
$ws_ARM.Range("N32").Value = -26753.357

# ARM row 61
$ws_ARM.Range("H61").Value = 983.5208
$ws_ARM.Range("I61").Value = 934.8205
$ws_ARM.Range("J61").Value = 1194.5555
$ws_ARM.Range("K61").Value = 934.8205
$ws_ARM.Range("L61").Value = 1194.5555
$ws_ARM.Range("M61").Value = -722.8205
$ws_ARM.Range("N61").Value = -1618.5555

# ARM row 74
$ws_ARM.Range("H74").Value = 2293.7446
$ws_ARM.Range("I74").Value = 2051.3125
$ws_ARM.Range("J74").Value = 2810.9333
$ws_ARM.Range("K74").Value = 2051.3125
$ws_ARM.Range("L74").Value = 2810.9333
$ws_ARM.Range("M74").Value = -1177.3125
$ws_ARM.Range("N74").Value = -4558.933300000001

# ARM row 77
$ws_ARM.Range("H77").Value = 2293.7446
$ws_ARM.Range("I77").Value = 2051.3125
$ws_ARM.Range("J77").Value = 2810.9333
$ws_ARM.Range("K77").Value = 10256.5625
$ws_ARM.Range("L77").Value = 14054.6665
$ws_ARM.Range("M77").Value = -5888.5625
$ws_ARM.Range("N77").Value = -22790.6665

# ARM row 123
$ws_ARM.Range("H123").Value = 60300
$ws_ARM.Range("J123").Value = 60300
$ws_ARM.Range("L123").Value = 60300
$ws_ARM.Range("N123").Value = -70100

# ARM row 132
$ws_ARM.Range("H132").Value = 3014.5642
$ws_ARM.Range("I132").Value = 2460.4707
$ws_ARM.Range("J132").Value = 3442.7273
$ws_ARM.Range("K132").Value = 7381.4121
$ws_ARM.Range("L132").Value = 10328.1819
$ws_ARM.Range("M132").Value = -4851.4121
$ws_ARM.Range("N132").Value = -15388.1819

# ARM row 136
$ws_ARM.Range("H136").Value = 983.5208
$ws_ARM.Range("I136").Value = 934.8205
$ws_ARM.Range("J136").Value = 1194.5555
$ws_ARM.Range("K136").Value = 2804.4615
$ws_ARM.Range("L136").Value = 3583.6665
$ws_ARM.Range("M136").Value = -254.4615000000003
$ws_ARM.Range("N136").Value = -8683.666499999999

# BSM row 86
$ws_BSM.Range("H86").Value = 8000809
$ws_BSM.Range("I86").Value = 10000676
$ws_BSM.Range("J86").Value = 1340
$ws_BSM.Range("K86").Value = 10000676
$ws_BSM.Range("L86").Value = 1340
$ws_BSM.Range("M86").Value = -9999553
$ws_BSM.Range("N86").Value = -3586

# BSM row 89
$ws_BSM.Range("H89").Value = 8000809
$ws_BSM.Range("I89").Value = 10000676
$ws_BSM.Range("J89").Value = 1340
$ws_BSM.Range("K89").Value = 50003380
$ws_BSM.Range("L89").Value = 6700
$ws_BSM.Range("M89").Value = -49997764
$ws_BSM.Range("N89").Value = -17932

# BSM row 105
$ws_BSM.Range("H105").Value = 2227.027
$ws_BSM.Range("I105").Value = 569.65515
$ws_BSM.Range("J105").Value = 8235
$ws_BSM.Range("K105").Value = 569.65515
$ws_BSM.Range("L105").Value = 8235
$ws_BSM.Range("M105").Value = 1177.34485
$ws_BSM.Range("N105").Value = -11729

# BSM row 107
$ws_BSM.Range("H107").Value = 839
$ws_BSM.Range("I107").Value = 590.2759
$ws_BSM.Range("J107").Value = 1494.7273
$ws_BSM.Range("K107").Value = 590.2759
$ws_BSM.Range("L107").Value = 1494.7273
$ws_BSM.Range("M107").Value = 1329.7241
$ws_BSM.Range("N107").Value = -5334.7273

# BSM row 134
$ws_BSM.Range("H134").Value = 1334.3541
$ws_BSM.Range("I134").Value = 1068.5
$ws_BSM.Range("J134").Value = 2663.625
$ws_BSM.Range("K134").Value = 3205.5
$ws_BSM.Range("L134").Value = 7990.875
$ws_BSM.Range("M134").Value = -670.5
$ws_BSM.Range("N134").Value = -13060.875

# CRP row 31
$ws_CRP.Range("H31").Value = 2514309.8
$ws_CRP.Range("I31").Value = 5154871.5
$ws_CRP.Range("J31").Value = 2555.756
$ws_CRP.Range("K31").Value = 5154871.5
$ws_CRP.Range("L31").Value = 2555.756
$ws_CRP.Range("M31").Value = -5154576.5
$ws_CRP.Range("N31").Value = -3145.756

# CRP row 34
$ws_CRP.Range("H34").Value = 2514309.8
$ws_CRP.Range("I34").Value = 5154871.5
$ws_CRP.Range("J34").Value = 2555.756
$ws_CRP.Range("K34").Value = 5154871.5
$ws_CRP.Range("L34").Value = 2555.756
$ws_CRP.Range("M34").Value = -5154669.5
$ws_CRP.Range("N34").Value = -2959.756

# CRP row 48
$ws_CRP.Range("H48").Value = 6142.143
$ws_CRP.Range("J48").Value = 6142.143
$ws_CRP.Range("L48").Value = 6142.143
$ws_CRP.Range("N48").Value = -7094.143

# CRP row 58
$ws_CRP.Range("H58").Value = 787.69696
$ws_CRP.Range("I58").Value = 419.9091
$ws_CRP.Range("J58").Value = 1523.2727
$ws_CRP.Range("K58").Value = 419.9091
$ws_CRP.Range("L58").Value = 1523.2727
$ws_CRP.Range("M58").Value = -216.9091
$ws_CRP.Range("N58").Value = -1929.2727

# CRP row 99
$ws_CRP.Range("H99").Value = 1431820.4
$ws_CRP.Range("I99").Value = 1882616.6
$ws_CRP.Range("J99").Value = 4299
$ws_CRP.Range("K99").Value = 1882616.6
$ws_CRP.Range("L99").Value = 4299
$ws_CRP.Range("M99").Value = -1881118.6
$ws_CRP.Range("N99").Value = -7295

# CRP row 126
$ws_CRP.Range("H126").Value = 1431820.4
$ws_CRP.Range("I126").Value = 1882616.6
$ws_CRP.Range("J126").Value = 4299
$ws_CRP.Range("K126").Value = 5647849.800000001
$ws_CRP.Range("L126").Value = 12897
$ws_CRP.Range("M126").Value = -5645379.800000001
$ws_CRP.Range("N126").Value = -17837

# CRP row 136
$ws_CRP.Range("H136").Value = 787.69696
$ws_CRP.Range("I136").Value = 419.9091
$ws_CRP.Range("J136").Value = 1523.2727
$ws_CRP.Range("K136").Value = 1259.7273
$ws_CRP.Range("L136").Value = 4569.8181
$ws_CRP.Range("M136").Value = 1290.2727
$ws_CRP.Range("N136").Value = -9669.8181

# CUL row 50
$ws_CUL.Range("H50").Value = 233.09091
$ws_CUL.Range("I50").Value = 48.333332
$ws_CUL.Range("J50").Value = 302.375
$ws_CUL.Range("K50").Value = 144.999996
$ws_CUL.Range("L50").Value = 907.125
$ws_CUL.Range("M50").Value = 336.000004
$ws_CUL.Range("N50").Value = -1869.125

# CUL row 53
$ws_CUL.Range("H53").Value = 233.09091
$ws_CUL.Range("I53").Value = 48.333332
$ws_CUL.Range("J53").Value = 302.375
$ws_CUL.Range("K53").Value = 144.999996
$ws_CUL.Range("L53").Value = 907.125
$ws_CUL.Range("M53").Value = 336.000004
$ws_CUL.Range("N53").Value = -1869.125

# GSM row 2
$ws_GSM.Range("H2").Value = 119.68421
$ws_GSM.Range("I2").Value = 104.3
$ws_GSM.Range("J2").Value = 136.77777
$ws_GSM.Range("K2").Value = 104.3
$ws_GSM.Range("L2").Value = 136.77777
$ws_GSM.Range("M2").Value = 8.700000000000003
$ws_GSM.Range("N2").Value = -362.77777

# GSM row 70
$ws_GSM.Range("H70").Value = 17395746
$ws_GSM.Range("I70").Value = 36367620
$ws_GSM.Range("J70").Value = 4862.3335
$ws_GSM.Range("K70").Value = 36367620
$ws_GSM.Range("L70").Value = 4862.3335
$ws_GSM.Range("M70").Value = -36367350
$ws_GSM.Range("N70").Value = -5402.3335

# GSM row 73
$ws_GSM.Range("H73").Value = 17395746
$ws_GSM.Range("I73").Value = 36367620
$ws_GSM.Range("J73").Value = 4862.3335
$ws_GSM.Range("K73").Value = 36367620
$ws_GSM.Range("L73").Value = 4862.3335
$ws_GSM.Range("M73").Value = -36366684
$ws_GSM.Range("N73").Value = -6734.3335

# GSM row 102
$ws_GSM.Range("H102").Value = 1561.6774
$ws_GSM.Range("I102").Value = 1353.2727
$ws_GSM.Range("K102").Value = 1353.2727
$ws_GSM.Range("M102").Value = 268.7273

# GSM row 132
$ws_GSM.Range("H132").Value = 2329.6924
$ws_GSM.Range("I132").Value = 1826
$ws_GSM.Range("J132").Value = 3463
$ws_GSM.Range("K132").Value = 5478
$ws_GSM.Range("L132").Value = 10389
$ws_GSM.Range("M132").Value = -2948
$ws_GSM.Range("N132").Value = -15449

# LTW row 40
$ws_LTW.Range("H40").Value = 1223.0476
$ws_LTW.Range("I40").Value = 1145.8823
$ws_LTW.Range("J40").Value = 1551
$ws_LTW.Range("K40").Value = 1145.8823
$ws_LTW.Range("L40").Value = 1551
$ws_LTW.Range("M40").Value = -1009.8823
$ws_LTW.Range("N40").Value = -1823

# LTW row 56
$ws_LTW.Range("H56").Value = 20000
$ws_LTW.Range("J56").Value = 20000
$ws_LTW.Range("L56").Value = 20000
$ws_LTW.Range("N56").Value = -21382

# LTW row 122
$ws_LTW.Range("H122").Value = 2566.4285
$ws_LTW.Range("I122").Value = 2566.4285
$ws_LTW.Range("J122").Value = 0
$ws_LTW.Range("K122").Value = 7699.2855
$ws_LTW.Range("L122").Value = 0
$ws_LTW.Range("M122").Value = ""
$ws_LTW.Range("N122").Value = -5249.2855

# LTW row 132
$ws_LTW.Range("H132").Value = 10781343
$ws_LTW.Range("I132").Value = 25009364
$ws_LTW.Range("J132").Value = 2539.2727
$ws_LTW.Range("K132").Value = 75028092
$ws_LTW.Range("L132").Value = 7617.8181
$ws_LTW.Range("M132").Value = -75025562
$ws_LTW.Range("N132").Value = -12677.8181

# WVR row 58
$ws_WVR.Range("H58").Value = 4264.1665
$ws_WVR.Range("I58").Value = 585
$ws_WVR.Range("J58").Value = 5000
$ws_WVR.Range("K58").Value = 585
$ws_WVR.Range("L58").Value = 5000
$ws_WVR.Range("M58").Value = -277
$ws_WVR.Range("N58").Value = -5616

# WVR row 61
$ws_WVR.Range("H61").Value = 7775.5
$ws_WVR.Range("I61").Value = 8051
$ws_WVR.Range("J61").Value = 7500
$ws_WVR.Range("K61").Value = 8051
$ws_WVR.Range("L61").Value = 7500
$ws_WVR.Range("M61").Value = -7759
$ws_WVR.Range("N61").Value = -8084

# WVR row 70
$ws_WVR.Range("H70").Value = 10000
$ws_WVR.Range("I70").Value = 10000
$ws_WVR.Range("J70").Value = 0
$ws_WVR.Range("K70").Value = 10000
$ws_WVR.Range("L70").Value = ""
$ws_WVR.Range("N70").Value = 0
$ws_WVR.Range("M70").Value = -9685

# WVR row 73
$ws_WVR.Range("H73").Value = 10000
$ws_WVR.Range("I73").Value = 10000
$ws_WVR.Range("J73").Value = 0
$ws_WVR.Range("K73").Value = 10000
$ws_WVR.Range("L73").Value = ""
$ws_WVR.Range("N73").Value = 0
$ws_WVR.Range("M73").Value = -8908

# WVR row 122
$ws_WVR.Range("H122").Value = 1350.1538
$ws_WVR.Range("I122").Value = 1350.1538
$ws_WVR.Range("J122").Value = 0
$ws_WVR.Range("K122").Value = 4050.4614
$ws_WVR.Range("L122").Value = 0
$ws_WVR.Range("M122").Value = ""
$ws_WVR.Range("N122").Value = -1600.4614

# WVR row 123
$ws_WVR.Range("H123").Value = 50765
$ws_WVR.Range("J123").Value = 50765
$ws_WVR.Range("L123").Value = 50765
$ws_WVR.Range("N123").Value = -60565

# WVR row 132
$ws_WVR.Range("H132").Value = 1016.6462
$ws_WVR.Range("I132").Value = 654.1739
$ws_WVR.Range("J132").Value = 1894.2106
$ws_WVR.Range("K132").Value = 1962.5217
$ws_WVR.Range("L132").Value = 5682.6318
$ws_WVR.Range("M132").Value = 567.4783
$ws_WVR.Range("N132").Value = -10742.6318

# WVR row 136
$ws_WVR.Range("H136").Value = 787.0714
$ws_WVR.Range("I136").Value = 447.625
$ws_WVR.Range("J136").Value = 1239.6666
$ws_WVR.Range("K136").Value = 1342.875
$ws_WVR.Range("L136").Value = 3718.9998
$ws_WVR.Range("M136").Value = 1207.125
$ws_WVR.Range("N136").Value = -8818.9998
